$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.760.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.33%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.117.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +11.19%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9995"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5190"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.84%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4362"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09040"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +8.45%  "
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.27%  "
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.180"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.95"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.115.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +11.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.795"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.672"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "97.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.95%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001142"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.38%  "
$ws.Range("B18").Value = "BinanceUSD"
$ws.Range("C18").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06620"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.06%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.447"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +8.85%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9995"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.964.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.365.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +11.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.276"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.570"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +12.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "164.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.76%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.192"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1069"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.249"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.534"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +29.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.892"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02588"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.620"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06791"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "9.537"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +11.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.68"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +11.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2252"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6799"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.255"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9991"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6323"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.256"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.665"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.277"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "127.85"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "83.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.11%  "
